$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.157.33'
$ws.Range("E2").Value = '  +0.93%  '

$ws.Range("D3").Value = '2.805.72'
$ws.Range("E3").Value = '  +1.61%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.68'
$ws.Range("E5").Value = '  +6.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '115.75'
$ws.Range("E6").Value = '  -1.72%  '

$ws.Range("E7").Value = '  +2.76%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.599'
$ws.Range("E9").Value = '  +3.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.50'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("E11").Value = '  +3.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.97'
$ws.Range("E12").Value = '  -1.41%  '

$ws.Range("E13").Value = '  +1.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.91'
$ws.Range("E14").Value = '  +3.66%  '

$ws.Range("D15").Value = '3.243.70'
$ws.Range("E15").Value = '  +1.61%  '

$ws.Range("D16").Value = '2.801.17'
$ws.Range("E16").Value = '  +0.92%  '

$ws.Range("E17").Value = '  +1.44%  '

$ws.Range("D18").Value = '52.106.28'
$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.32'
$ws.Range("E19").Value = '  +6.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.16'
$ws.Range("E20").Value = '  +5.72%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.57'
$ws.Range("E21").Value = '  -1.01%  '

$ws.Range("D22").Value = '0.0₃0983'
$ws.Range("E22").Value = '  +1.84%  '

$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.85'
$ws.Range("E24").Value = '  -2.70%  '

$ws.Range("E25").Value = '  +4.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.95'
$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.30'
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.26'
$ws.Range("E29").Value = '  +1.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.140'
$ws.Range("E30").Value = '  -0.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.52'
$ws.Range("E31").Value = '  -3.11%  '

$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.46'
$ws.Range("E32").Value = '  +0.02%  '

$ws.Range("B33").Value = 'VeChain'
$ws.Range("C33").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0449'
$ws.Range("E33").Value = '  +30.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.80'
$ws.Range("E34").Value = '  +3.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0835'
$ws.Range("E35").Value = '  -0.26%  '

$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.15%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.11'
$ws.Range("E37").Value = '  -0.21%  '

$ws.Range("E38").Value = '  -0.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.78'
$ws.Range("E39").Value = '  -3.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.25'
$ws.Range("E40").Value = '  -0.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.60'
$ws.Range("E41").Value = '  +9.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.66'
$ws.Range("E42").Value = '  +0.49%  '

$ws.Range("E43").Value = '  +2.11%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.31'
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '125.39'
$ws.Range("E45").Value = '  -4.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.34'
$ws.Range("E46").Value = '  -0.73%  '

$ws.Range("D47").Value = '2.057.40'
$ws.Range("E47").Value = '  -2.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.31'
$ws.Range("E48").Value = '  +1.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.950'
$ws.Range("E49").Value = '  +8.93%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.58'
$ws.Range("E50").Value = '  -0.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.15'
$ws.Range("E51").Value = '  +1.55%  '
